$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data edits on the "requirements" sheet/table ---

# B3: 1 -> 5
$ws.Range("B3").Value = 5

# H3: numeric 1 -> text "1" (new text-formatted style w/ same border as before)
$ws.Range("H3").NumberFormat = "@"
$ws.Range("H3").Value = "1"

# B4: numeric 1 -> text "tet" (new shared string)
$ws.Range("B4").Value = "tet"

# F5: 1 -> 5.7
$ws.Range("F5").Value = 5.7

# --- Remove the data validation rules (finished validation work, no longer needed) ---
[void]$ws.Cells.Validation.Delete()

# --- Update the active selection shown when the sheet is opened ---
[void]$ws.Range("J20").Select()
